$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(
    @{ r=404; A=27340; B=45457; C=4; D="Minnesota ROKKR";       E="Vista";    F="Sub Base"; G="Invasion"; H="6 Star";   I="Highrise" },
    @{ r=405; A=27340; B=45457; C=4; D="New York Subliners";    E="6 Star";   F="Karachi";  G="Rio";      H="Highrise"; J="Karachi" },
    @{ r=406; A=27341; B=45457; C=4; D="Carolina Royal Ravens"; E="6 Star";   F="Vista";    G="6 Star";   H="Karachi";  I="Invasion" },
    @{ r=407; A=27341; B=45457; C=4; D="OpTic Texas";           E="Karachi";  F="Rio";      G="Rio";      H="Highrise"; J="Highrise" },
    @{ r=408; A=27342; B=45457; C=4; D="Boston Breach";         E="Rio";      F="Karachi";  G="6 Star";   H="Invasion"; I="Highrise" },
    @{ r=409; A=27342; B=45457; C=4; D="Atlanta FaZe";          E="Sub Base"; F="Vista";    G="Highrise"; H="Karachi";  J="Karachi" }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    if ($row.ContainsKey("I")) { $ws.Cells.Item($r, 9).Value = $row.I }
    if ($row.ContainsKey("J")) { $ws.Cells.Item($r, 10).Value = $row.J }
}

# Update the view: frozen pane scroll position and active selection, matching the
# post-edit state where the user scrolled near the bottom of the newly-added data.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 396
$ws.Range("F413").Select() | Out-Null

$wb.Save()
